# summarize findings for twosigma, uber
# Adds a new "Collections" category (RQ1 table + breakdown rows) and two new
# RQ2 "(Collections) ..." fix reasons, bumps several counts, and adds a new
# "Clojure" row to the Programming Language table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New label cells -------------------------------------------------
# Order matters: it controls the order new entries are appended to the
# workbook's shared-string table (Clojure, Collections, then the two new
# "(Collections) ..." fix-reason strings).
$ws.Range("B29").Value = "Clojure"
$ws.Range("B11").Value = "Collections"
$ws.Range("E16").Value = "(Collections) sort / establish ordering"
$ws.Range("E17").Value = "(Collections) limit elements"

# --- Updated / new counts --------------------------------------------
$ws.Range("C4").Value = 12    # Tolerance
$ws.Range("F4").Value = 7     # (Tolerance) Increase acceptance levels in assert statements
$ws.Range("C8").Value = 5     # Concurrency
$ws.Range("C11").Value = 2    # Collections (new)
$ws.Range("F12").Value = 5    # (Concurrency) locks
$ws.Range("F16").Value = 1    # (Collections) sort / establish ordering (new)
$ws.Range("F17").Value = 1    # (Collections) limit elements (new)
$ws.Range("C22").Value = 5    # Python
$ws.Range("C27").Value = 4    # Go
$ws.Range("C29").Value = 4    # Clojure (new)

# --- Grow the RQ2 ("Table13") table to include the two new rows ------
$lo = $ws.ListObjects.Item("Table13")
$lo.Resize($ws.Range("E3:F17"))

# --- View state: scroll up and move the active selection -------------
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 1
$ws.Range("F4").Select()
